$d = $word.ActiveDocument

# --- helpers -------------------------------------------------------------

function Get-EnclosingParagraph($pos) {
    foreach ($pp in $d.Paragraphs) {
        if ($pp.Range.Start -le $pos -and $pp.Range.End -gt $pos) {
            return $pp
        }
    }
    return $null
}

# Inspects the real OOXML of the paragraph containing $pos to see whether the
# run structure right before $pos is an empty (text-less) run, e.g. <w:r/>.
# Find/Replace on a run's text tends to merge/drop such adjacent empty runs,
# so we detect their presence beforehand and restore them afterward.
function Test-EmptyRunBefore($pos) {
    $para = Get-EnclosingParagraph $pos
    if ($para -eq $null) { return $false }

    $xml = $para.Range.WordOpenXML
    if ($xml -notmatch '(?s)<w:body>(.*)</w:body>') { return $false }
    $body = $Matches[1]
    if ($body -notmatch '(?s)<w:p[ >].*?</w:p>') { return $false }
    $pxml = $Matches[0]

    # strip <w:pPr>...</w:pPr> so we only look at the run sequence
    $noPPr = [System.Text.RegularExpressions.Regex]::Replace($pxml, '(?s)<w:pPr>.*?</w:pPr>', '')

    # true when the paragraph's run list begins with a run that carries no <w:t>,
    # i.e. a self-closing <w:r/> or an explicit empty <w:r></w:r> (optionally
    # with run properties but never text)
    if ([System.Text.RegularExpressions.Regex]::IsMatch($noPPr, '(?s)^<w:p\b[^>]*>\s*<w:r\b[^>]*/>')) {
        return $true
    }
    if ([System.Text.RegularExpressions.Regex]::IsMatch($noPPr, '(?s)^<w:p\b[^>]*>\s*<w:r\b[^>]*>\s*</w:r>')) {
        return $true
    }
    if ([System.Text.RegularExpressions.Regex]::IsMatch($noPPr, '(?s)^<w:p\b[^>]*>\s*<w:r\b[^>]*>\s*<w:rPr>.*?</w:rPr>\s*</w:r>')) {
        return $true
    }
    return $false
}

function Insert-EmptyRun($pos) {
    $r = $d.Range($pos, $pos)
    $xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p><w:r/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($xml)
}

# Replaces every occurrence of $oldText with $newText (plain-text match,
# formatting of each run is preserved by Word's own Find/Replace) and heals
# any empty run that Find/Replace drops from immediately before the match.
function Replace-AllText($oldText, $newText) {
    $searchStart = 0
    while ($true) {
        $full = $d.Content.Text
        $idx = $full.IndexOf($oldText, $searchStart)
        if ($idx -lt 0) { break }
        $end = $idx + $oldText.Length

        $hadEmptyRun = Test-EmptyRunBefore $idx

        $r = $d.Range($idx, $end)
        $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 1) | Out-Null

        if ($hadEmptyRun) {
            $full2 = $d.Content.Text
            $newIdx = $full2.IndexOf($newText, $idx)
            if ($newIdx -ge 0 -and -not (Test-EmptyRunBefore $newIdx)) {
                Insert-EmptyRun $newIdx
            }
        }

        $searchStart = $idx + $newText.Length
    }
}

# --- edits -----------------------------------------------------------------

Replace-AllText "Play Golden Clover Free - Review & Features | Onlyplay" "Play Golden Clover for Free - Review & Gameplay"
Replace-AllText "Innovative lottery-style mechanics" "Innovative and simple lottery-style mechanics"
Replace-AllText "High RTP of 96.2%" "Relatively high RTP of 96.2%"
Replace-AllText "Interesting payouts for simplicity lovers" "Interesting payouts for players who love simplicity"
Replace-AllText "Perfect for players seeking new mechanism" "Ideal for players looking for a new mechanism"
Replace-AllText "Cartoon-style graphics might not suit everyone" "Cartoon-style graphics might not appeal to players who love sophisticated design"
Replace-AllText "Play Golden Clover for free and read our unbiased review about its features and gameplay. Developed by Onlyplay, it is an innovative and simple slot game." "Read our review of Golden Clover, a slot game with innovative lottery-style mechanics. Play for free!"
